$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R12").Value = 5

$ws.Range("R16").Value = 8
$ws.Range("T16").Value = 0
$ws.Range("U16").Value = 3

$ws.Range("R18").Value = 2
$ws.Range("T18").Value = 0

$ws.Range("R21").Value = 4
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 4

$ws.Range("R28").Value = 2

$ws.Range("R32").Value = 1

$ws.Range("R35").Value = 4

$ws.Range("R36").Value = 6

$ws.Range("R39").Value = 1

$ws.Range("R42").Value = 1

$ws.Range("R43").Value = 6

$ws.Range("R45").Value = 2
$ws.Range("T45").Value = 0

$ws.Range("R46").Value = 1
$ws.Range("T46").Value = 0

$ws.Range("R54").Value = 1

$ws.Range("R57").Value = 1

$ws.Range("R59").Value = 2
$ws.Range("T59").Value = 0

$ws.Range("R69").Value = 1
$ws.Range("T69").Value = 2

$ws.Range("R70").Value = 1

$ws.Range("R73").Value = 6
$ws.Range("T73").Value = 3
$ws.Range("U73").Value = 13

$ws.Range("R79").Value = 1

$ws.Range("R82").Value = 1
$ws.Range("T82").Value = 1

$ws.Range("R84").Value = 2

$ws.Range("R86").Value = 1

$ws.Range("R88").Value = 25
$ws.Range("T88").Value = 0
$ws.Range("U88").Value = 25

$ws.Range("R89").Value = 43
$ws.Range("T89").Value = 0

$ws.Range("R90").Value = 2

$ws.Range("R91").Value = 5
$ws.Range("T91").Value = 0

$ws.Range("R92").Value = 10
$ws.Range("T92").Value = 0

$ws.Range("C97").Value = 267
